# Scheduled-runner style market-price refresh across the per-job leve sheets.
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N) for the
# rows whose item market data changed since the last snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 33340166
$ws.Range("I48").Value = 7000
$ws.Range("J48").Value = 50006750
$ws.Range("K48").Value = 21000
$ws.Range("L48").Value = 150020250
$ws.Range("M48").Value = -20708
$ws.Range("N48").Value = -150020834

$ws.Range("H56").Value = 33340166
$ws.Range("I56").Value = 7000
$ws.Range("J56").Value = 50006750
$ws.Range("K56").Value = 21000
$ws.Range("L56").Value = 150020250
$ws.Range("M56").Value = -20466
$ws.Range("N56").Value = -150021318

$ws.Range("H80").Value = 66674500
$ws.Range("I80").Value = 200000670
$ws.Range("J80").Value = 11409.5
$ws.Range("K80").Value = 600002010
$ws.Range("L80").Value = 34228.5
$ws.Range("M80").Value = -600001012
$ws.Range("N80").Value = -36224.5

$ws.Range("H83").Value = 66674500
$ws.Range("I83").Value = 200000670
$ws.Range("J83").Value = 11409.5
$ws.Range("K83").Value = 1800006030
$ws.Range("L83").Value = 102685.5
$ws.Range("M83").Value = -1800001038
$ws.Range("N83").Value = -112669.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34193910
$ws.Range("I32").Value = 37880092
$ws.Range("J32").Value = 14288536
$ws.Range("K32").Value = 37880092
$ws.Range("L32").Value = 14288536
$ws.Range("M32").Value = -37879805

$ws.Range("H43").Value = 12138.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12138.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12138.5
$ws.Range("N43").Value = -12764.5
$ws.Range("M43").ClearContents()

$ws.Range("H44").Value = 69780.2
$ws.Range("I44").Value = 69000
$ws.Range("J44").Value = 69975.25
$ws.Range("K44").Value = 69000
$ws.Range("L44").Value = 69975.25
$ws.Range("M44").Value = -68512
$ws.Range("N44").Value = -70951.25

$ws.Range("H45").Value = 4332
$ws.Range("I45").Value = 4651.9165
$ws.Range("J45").Value = 4165.087
$ws.Range("K45").Value = 4651.9165
$ws.Range("L45").Value = 4165.087
$ws.Range("M45").Value = -4274.9165
$ws.Range("N45").Value = -4919.087

$ws.Range("H61").Value = 4076.7666
$ws.Range("I61").Value = 3623.4614
$ws.Range("J61").Value = 4423.4116
$ws.Range("K61").Value = 3623.4614
$ws.Range("L61").Value = 4423.4116
$ws.Range("M61").Value = -3411.4614
$ws.Range("N61").Value = -4847.4116

$ws.Range("H110").Value = 2819
$ws.Range("I110").Value = 2038.2
$ws.Range("J110").Value = 3599.8
$ws.Range("K110").Value = 2038.2
$ws.Range("L110").Value = 3599.8
$ws.Range("M110").Value = 6.799999999999955
$ws.Range("N110").Value = -7689.8

$ws.Range("H114").Value = 59329.668
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 59329.668
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 59329.668
$ws.Range("N114").Value = -68007.66800000001

$ws.Range("H122").Value = 5550.5
$ws.Range("I122").Value = 5254.3335
$ws.Range("J122").Value = 5728.2
$ws.Range("K122").Value = 15763.0005
$ws.Range("L122").Value = 17184.6
$ws.Range("M122").Value = -13313.0005
$ws.Range("N122").Value = -22084.6

$ws.Range("H131").Value = 112480.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 112480.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 112480.5
$ws.Range("N131").Value = -122560.5

$ws.Range("H132").Value = 306092.2
$ws.Range("I132").Value = 336365.06
$ws.Range("J132").Value = 3363.3333
$ws.Range("K132").Value = 1009095.18
$ws.Range("L132").Value = 10089.9999
$ws.Range("M132").Value = -1006565.18

$ws.Range("H136").Value = 4076.7666
$ws.Range("I136").Value = 3623.4614
$ws.Range("J136").Value = 4423.4116
$ws.Range("K136").Value = 10870.3842
$ws.Range("L136").Value = 13270.2348
$ws.Range("M136").Value = -8320.3842
$ws.Range("N136").Value = -18370.2348

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2378.4614
$ws.Range("I99").Value = 1586.6666
$ws.Range("J99").Value = 3057.1428
$ws.Range("K99").Value = 1586.6666
$ws.Range("L99").Value = 3057.1428
$ws.Range("M99").Value = -88.66660000000002

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3982.4814
$ws.Range("I31").Value = 1586
$ws.Range("J31").Value = 5760.516
$ws.Range("K31").Value = 1586
$ws.Range("L31").Value = 5760.516
$ws.Range("M31").Value = -1291
$ws.Range("N31").Value = -6350.516

$ws.Range("H34").Value = 3982.4814
$ws.Range("I34").Value = 1586
$ws.Range("J34").Value = 5760.516
$ws.Range("K34").Value = 1586
$ws.Range("L34").Value = 5760.516
$ws.Range("M34").Value = -1384
$ws.Range("N34").Value = -6164.516

$ws.Range("H70").Value = 45666.668
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 45666.668
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 45666.668
$ws.Range("N70").Value = -46296.668

$ws.Range("H73").Value = 45666.668
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 45666.668
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 45666.668
$ws.Range("N73").Value = -47850.668

$ws.Range("H105").Value = 1795.3077
$ws.Range("I105").Value = 1795.3077
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1795.3077
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -48.30770000000007

$ws.Range("H122").Value = 3236.7144
$ws.Range("I122").Value = 2765.7856
$ws.Range("J122").Value = 3707.6428
$ws.Range("K122").Value = 8297.356800000001
$ws.Range("L122").Value = 11122.9284
$ws.Range("M122").Value = -5847.356800000001
$ws.Range("N122").Value = -16022.9284

$ws.Range("H125").Value = 43889.75
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 43889.75
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 43889.75
$ws.Range("N125").Value = -48809.75

$ws.Range("H132").Value = 2084.3914
$ws.Range("I132").Value = 2065.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 6196.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3666.5
$ws.Range("N132").Value = -12560

$ws.Range("H134").Value = 2232.8333
$ws.Range("I134").Value = 2279.6
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 6838.799999999999
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -4303.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 199
$ws.Range("I50").Value = 38.8
$ws.Range("J50").Value = 1000
$ws.Range("K50").Value = 116.4
$ws.Range("L50").Value = 3000
$ws.Range("M50").Value = 364.6
$ws.Range("N50").Value = -3962

$ws.Range("H53").Value = 199
$ws.Range("I53").Value = 38.8
$ws.Range("J53").Value = 1000
$ws.Range("K53").Value = 116.4
$ws.Range("L53").Value = 3000
$ws.Range("M53").Value = 364.6
$ws.Range("N53").Value = -3962

$ws.Range("H86").Value = 113.5
$ws.Range("I86").Value = 113.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 340.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 845.5
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 113.5
$ws.Range("I89").Value = 113.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1021.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 4906.5
$ws.Range("N89").ClearContents()

$ws.Range("H116").Value = 432
$ws.Range("I116").Value = 404.8
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 1214.4
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 2227.6

$ws.Range("H131").Value = 1727.6666
$ws.Range("I131").Value = 1403.6
$ws.Range("J131").Value = 1775.3235
$ws.Range("K131").Value = 4210.799999999999
$ws.Range("L131").Value = 5325.970499999999
$ws.Range("M131").Value = 829.2000000000007
$ws.Range("N131").Value = -15405.9705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4199.6
$ws.Range("I122").Value = 2999.6667
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 8999.000100000001
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -6549.000100000001
$ws.Range("N122").Value = -22898.5

$ws.Range("H126").Value = 2950.842
$ws.Range("I126").Value = 2583.3572
$ws.Range("J126").Value = 3979.8
$ws.Range("K126").Value = 7750.071599999999
$ws.Range("L126").Value = 11939.4
$ws.Range("M126").Value = -5280.071599999999
$ws.Range("N126").Value = -16879.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2049.25
$ws.Range("I40").Value = 2049.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2049.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1913.25

$ws.Range("H125").Value = 78770
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 78770
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 78770
$ws.Range("N125").Value = -88610

$ws.Range("H132").Value = 716247
$ws.Range("I132").Value = 835288.2
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2505864.6
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2503334.6

$ws.Range("H136").Value = 8349.700000000001
$ws.Range("I136").Value = 2250
$ws.Range("J136").Value = 17499.25
$ws.Range("K136").Value = 6750
$ws.Range("L136").Value = 52497.75
$ws.Range("M136").Value = -4200
$ws.Range("N136").Value = -57597.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 20007
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 20007
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 20007
$ws.Range("N15").Value = -20583
$ws.Range("M15").ClearContents()

$ws.Range("H81").Value = 48519.41
$ws.Range("I81").Value = 78749.08
$ws.Range("J81").Value = 4854.3335
$ws.Range("K81").Value = 157498.16
$ws.Range("L81").Value = 9708.666999999999
$ws.Range("M81").Value = -156437.16
$ws.Range("N81").Value = -11830.667

$ws.Range("H84").Value = 48519.41
$ws.Range("I84").Value = 78749.08
$ws.Range("J84").Value = 4854.3335
$ws.Range("K84").Value = 787490.8
$ws.Range("L84").Value = 48543.335
$ws.Range("M84").Value = -782186.8
$ws.Range("N84").Value = -59151.335

$ws.Range("H107").Value = 975
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 975
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2925
$ws.Range("N107").Value = -6765
$ws.Range("M107").ClearContents()

$ws.Range("H130").Value = 61717
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 61717
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 61717
$ws.Range("N130").Value = -71757

$ws.Range("H132").Value = 38987.285
$ws.Range("I132").Value = 44452.043
$ws.Range("J132").Value = 6198.75
$ws.Range("K132").Value = 133356.129
$ws.Range("L132").Value = 18596.25
$ws.Range("M132").Value = -130826.129

$ws.Range("H136").Value = 2214.7188
$ws.Range("I136").Value = 1354.35
$ws.Range("J136").Value = 3648.6667
$ws.Range("K136").Value = 4063.05
$ws.Range("L136").Value = 10946.0001
$ws.Range("M136").Value = -1513.05
$ws.Range("N136").Value = -16046.0001
